$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2:C89").Value = 0
